# Auto-generated Excel COM-interop edit script
# Applies the "Horarios actualizados Linea 141 - 114" scrape update
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 22:06:14"
$ws1.Range("A3").Value = "Total filas: 161"

$ws1.Cells.Item(44, 1).Value = "17:35:09"
$ws1.Cells.Item(44, 2).Value = "18:03"
$ws1.Cells.Item(44, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(44, 4).Value = 28
$ws1.Cells.Item(44, 5).Value = "LP1912"

$ws1.Cells.Item(45, 1).Value = "16:52:37"
$ws1.Cells.Item(45, 2).Value = "18:03"
$ws1.Cells.Item(45, 3).Value = "17_ROMERO"
$ws1.Cells.Item(45, 4).Value = 71
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(85, 1).Value = "17:54:43"
$ws1.Cells.Item(85, 2).Value = "19:30"
$ws1.Cells.Item(85, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(85, 4).Value = 96
$ws1.Cells.Item(85, 5).Value = "LP1912"

$ws1.Cells.Item(86, 1).Value = "17:54:43"
$ws1.Cells.Item(86, 2).Value = "19:30"
$ws1.Cells.Item(86, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(86, 4).Value = 96
$ws1.Cells.Item(86, 5).Value = "LP1912"

$ws1.Cells.Item(105, 1).Value = "19:47:58"
$ws1.Cells.Item(105, 2).Value = "20:00"
$ws1.Cells.Item(105, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(105, 4).Value = 13
$ws1.Cells.Item(105, 5).Value = "LP1912"

$ws1.Cells.Item(106, 1).Value = "18:10:41"
$ws1.Cells.Item(106, 2).Value = "20:00"
$ws1.Cells.Item(106, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(106, 4).Value = 110
$ws1.Cells.Item(106, 5).Value = "LP1912"

$ws1.Cells.Item(125, 1).Value = "18:52:04"
$ws1.Cells.Item(125, 2).Value = "20:44"
$ws1.Cells.Item(125, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(125, 4).Value = 112
$ws1.Cells.Item(125, 5).Value = "LP1912"

$ws1.Cells.Item(126, 1).Value = "19:11:59"
$ws1.Cells.Item(126, 2).Value = "20:44"
$ws1.Cells.Item(126, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(126, 4).Value = 93
$ws1.Cells.Item(126, 5).Value = "LP1912"

$ws1.Cells.Item(154, 1).Value = "22:06:14"
$ws1.Cells.Item(154, 2).Value = "22:21"
$ws1.Cells.Item(154, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(154, 4).Value = 15
$ws1.Cells.Item(154, 5).Value = "LP1912"

$ws1.Cells.Item(155, 1).Value = "22:06:14"
$ws1.Cells.Item(155, 2).Value = "22:22"
$ws1.Cells.Item(155, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(155, 4).Value = 16
$ws1.Cells.Item(155, 5).Value = "LP1912"

$ws1.Cells.Item(156, 1).Value = "20:45:46"
$ws1.Cells.Item(156, 2).Value = "22:25"
$ws1.Cells.Item(156, 3).Value = "15_ABASTO"
$ws1.Cells.Item(156, 4).Value = 100
$ws1.Cells.Item(156, 5).Value = "LP1912"

$ws1.Cells.Item(157, 1).Value = "20:32:02"
$ws1.Cells.Item(157, 2).Value = "22:26"
$ws1.Cells.Item(157, 3).Value = "15_ABASTO"
$ws1.Cells.Item(157, 4).Value = 114
$ws1.Cells.Item(157, 5).Value = "LP1912"

$ws1.Cells.Item(158, 1).Value = "20:45:46"
$ws1.Cells.Item(158, 2).Value = "22:30"
$ws1.Cells.Item(158, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(158, 4).Value = 105
$ws1.Cells.Item(158, 5).Value = "LP1912"

$ws1.Cells.Item(159, 1).Value = "20:32:02"
$ws1.Cells.Item(159, 2).Value = "22:31"
$ws1.Cells.Item(159, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(159, 4).Value = 119
$ws1.Cells.Item(159, 5).Value = "LP1912"

$ws1.Cells.Item(160, 1).Value = "22:06:14"
$ws1.Cells.Item(160, 2).Value = "22:41"
$ws1.Cells.Item(160, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(160, 4).Value = 35
$ws1.Cells.Item(160, 5).Value = "LP1912"

$ws1.Cells.Item(161, 1).Value = "20:52:24"
$ws1.Cells.Item(161, 2).Value = "22:49"
$ws1.Cells.Item(161, 3).Value = "14_ABASTO"
$ws1.Cells.Item(161, 4).Value = 117
$ws1.Cells.Item(161, 5).Value = "LP1912"

$ws1.Cells.Item(162, 1).Value = "22:06:14"
$ws1.Cells.Item(162, 2).Value = "23:06"
$ws1.Cells.Item(162, 3).Value = "15_ABASTO"
$ws1.Cells.Item(162, 4).Value = 60
$ws1.Cells.Item(162, 5).Value = "LP1912"

$ws1.Cells.Item(163, 1).Value = "22:06:14"
$ws1.Cells.Item(163, 2).Value = "23:10"
$ws1.Cells.Item(163, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(163, 4).Value = 64
$ws1.Cells.Item(163, 5).Value = "LP1912"

$ws1.Cells.Item(164, 1).Value = "22:06:14"
$ws1.Cells.Item(164, 2).Value = "23:19"
$ws1.Cells.Item(164, 3).Value = "14_ABASTO"
$ws1.Cells.Item(164, 4).Value = 73
$ws1.Cells.Item(164, 5).Value = "LP1912"

$ws1.Cells.Item(165, 1).Value = "22:06:14"
$ws1.Cells.Item(165, 2).Value = "23:44"
$ws1.Cells.Item(165, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(165, 4).Value = 98
$ws1.Cells.Item(165, 5).Value = "LP1912"

$ws1.Cells.Item(166, 1).Value = "22:06:14"
$ws1.Cells.Item(166, 2).Value = "23:49"
$ws1.Cells.Item(166, 3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(166, 4).Value = 103
$ws1.Cells.Item(166, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 22:06:14"
$ws2.Range("A3").Value = "Total filas: 20"

$ws2.Cells.Item(25, 1).Value = "22:06:14"
$ws2.Cells.Item(25, 2).Value = "23:44"
$ws2.Cells.Item(25, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(25, 4).Value = 98
$ws2.Cells.Item(25, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 22:06:14"
$ws3.Range("A3").Value = "Total filas: 20"

$ws3.Cells.Item(21, 1).Value = "22:06:14"
$ws3.Cells.Item(21, 2).Value = "22:08"
$ws3.Cells.Item(21, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(21, 4).Value = 2
$ws3.Cells.Item(21, 5).Value = "L6203"

$ws3.Cells.Item(22, 1).Value = "20:45:46"
$ws3.Cells.Item(22, 2).Value = "22:12"
$ws3.Cells.Item(22, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(22, 4).Value = 87
$ws3.Cells.Item(22, 5).Value = "L6173"

$ws3.Cells.Item(23, 1).Value = "20:32:02"
$ws3.Cells.Item(23, 2).Value = "22:13"
$ws3.Cells.Item(23, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(23, 4).Value = 101
$ws3.Cells.Item(23, 5).Value = "L6173"

$ws3.Cells.Item(24, 1).Value = "20:32:02"
$ws3.Cells.Item(24, 2).Value = "22:19"
$ws3.Cells.Item(24, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(24, 4).Value = 107
$ws3.Cells.Item(24, 5).Value = "L6173"

$ws3.Cells.Item(25, 1).Value = "22:06:14"
$ws3.Cells.Item(25, 2).Value = "22:23"
$ws3.Cells.Item(25, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(25, 4).Value = 17
$ws3.Cells.Item(25, 5).Value = "L6173"
